$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $found = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $found) {
        Write-Output "WARNING: text not found -> $find"
    }
}

# "Door een handtekening..." paragraph: drop "verder " and "gemeld" -> "genoemd"
Replace-Text "zaken die verder in het samenwerkingscontract gemeld worden" "zaken die in het samenwerkingscontract genoemd worden"

# "Roulatiesysteem: De notulist..." -> lowercase the "De"
Replace-Text "Roulatiesysteem: De notulist" "Roulatiesysteem: de notulist"

# "Bij het committen..." bullet rewrite
Replace-Text "van wat er veranderd is" "van de toevoeging, verandering of verwijdering"
Replace-Text "kan er makkelijk nagegaan" "kan er gemakkelijk nagegaan"
Replace-Text "fouten makkelijker gevonden worden" "fouten makkelijker verholpen worden"

# Add missing period after "... aanwezig zijn"
Replace-Text "Niet op tijd aanwezig = meer dan 5 minuten te laat aanwezig zijn" "Niet op tijd aanwezig = meer dan 5 minuten te laat aanwezig zijn."

# "wordt bij gehouden." -> "wordt bijgehouden."
Replace-Text "Afwezigheid / te laat komen wordt bij gehouden." "Afwezigheid / te laat komen wordt bijgehouden."

# Add a default footer with the date, right-aligned, using the "Voettekst" (Footer) style.
$sec = $d.Sections.First
$footer = $sec.Footers.Item(1)
$footer.Range.InsertAfter("Datum: 30 november 2010")
$footer.Range.Paragraphs.Item(1).Style = "Voettekst"
$footer.Range.Paragraphs.Item(1).Alignment = 2

Write-Output "Edit complete"
